# Weekly fruit/hortaliza data update:
# A new weekly price observation is inserted as row 263 (pushing the
# existing rows 263-362 down to 264-363), matching the canonical OOXML
# diff which adds one new data row to the "Poroto verde" price series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 263, shifting everything below it down.
$ws.Rows.Item(263).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A263").Value = 8
$ws.Range("B263").Value = "Terminal La Palmera de La Serena"
$ws.Range("C263").Value = "Coquimbo"
$ws.Range("D263").Value = 45009
$ws.Range("E263").Value = 4
$ws.Range("F263").Value = 100112031
$ws.Range("G263").Value = "Poroto verde"
$ws.Range("H263").Value = "Magnum"
$ws.Range("I263").Value = "Primera"
$ws.Range("J263").Value = 400
$ws.Range("K263").Value = 23000
$ws.Range("L263").Value = 24000
$ws.Range("M263").Value = 23500
$ws.Range("N263").Value = "$/malla 25 kilos"
$ws.Range("O263").Value = "Provincia de Limarí"
$ws.Range("P263").Value = 940
$ws.Range("Q263").Value = 25
$ws.Range("R263").Value = "Hortaliza"
